# Atualização relatorio de entrega
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: only the "Alunos" column changes ---
$ws.Range("A4").Value = "Gabriel"

# --- Row 5: Alunos + Atividade change ---
$ws.Range("A5").Value = "Rafael"
$ws.Range("B5").Value = "Seleção da arquitetura, tecnologias e padrões"

# --- Row 6: Alunos + Atividade change ---
$ws.Range("A6").Value = "Gabriel, Willian"
$ws.Range("B6").Value = "Pesquisa sobre a arquitetura"

# --- Row 7: Atividade + Importância change; underline formatting moves away ---
$ws.Range("B7").Value = "Pesquisa sobre as tecnologias e padrões"
$ws.Range("C7").Value = "Alta"
$ws.Range("A7").Font.Underline = -4142

# --- Row 8: fill in the previously-blank placeholder row ---
$ws.Range("A8").Value = "Gabriel, Willian, Rafael"
$ws.Range("B8").Value = "Justificativa do uso da arquitetura, tecnologias e padrões"
$ws.Range("C8").Value = "Alta"

# --- Row 9: brand-new row appended at the bottom ---
# Copy the formatting from row 4 (a normal data row) down into row 9 first,
# then overwrite the values - this preserves the exact cell styles (incl. the
# "EFEFEF"-less normal style used by columns B/C).
$ws.Range("A4:C4").Copy($ws.Range("A9:C9"))

$ws.Range("A9").Value = "Willian"
$ws.Range("B9").Value = "Formalização do documento de entrega"
$ws.Range("C9").Value = "Média"

$ws.Rows.Item(9).RowHeight = 15.75

# Row 9, column A is the new "last row" of the table, so it gets the
# underline styling that used to sit on row 7 / column A.
$ws.Range("A9").Font.Name = "Arial"
$ws.Range("A9").Font.Size = 10
$ws.Range("A9").Font.Color = 0
$ws.Range("A9").Font.Underline = 2

# --- move the active selection, matching the saved cursor position ---
$ws.Range("A11").Select()
